$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = -0.4421614382907215
$ws.Range("J2").Value = 0.2324881689575329
$ws.Range("K2").Value = -0.1758749011479283
$ws.Range("L2").Value = 2.716096556723631
